# Applies updated crypto price/volume/hour data to Sheet1 (rows 2-51)
# Values are stored as text (matching the source sheet's inlineStr cells),
# so each write uses a leading apostrophe to force text entry, then clears
# the resulting "quote prefix" number format back to the default style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'309.32"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.03%"
$ws.Range("E2").ClearFormats()
$ws.Range("G2").Value = "'12"
$ws.Range("G2").ClearFormats()

$ws.Range("D3").Value = "'41.31"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'0.70%"
$ws.Range("E3").ClearFormats()
$ws.Range("G3").Value = "'12"
$ws.Range("G3").ClearFormats()

$ws.Range("D4").Value = "'5.198"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'1.62%"
$ws.Range("E4").ClearFormats()
$ws.Range("G4").Value = "'12"
$ws.Range("G4").ClearFormats()

$ws.Range("D5").Value = "'0.07709"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'0.94%"
$ws.Range("E5").ClearFormats()
$ws.Range("G5").Value = "'12"
$ws.Range("G5").ClearFormats()

$ws.Range("D6").Value = "'1.647"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'2.61%"
$ws.Range("E6").ClearFormats()
$ws.Range("G6").Value = "'12"
$ws.Range("G6").ClearFormats()

$ws.Range("D7").Value = "'0.9156"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'1.29%"
$ws.Range("E7").ClearFormats()
$ws.Range("G7").Value = "'12"
$ws.Range("G7").ClearFormats()

$ws.Range("D8").Value = "'2.426"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'-3.60%"
$ws.Range("E8").ClearFormats()
$ws.Range("G8").Value = "'12"
$ws.Range("G8").ClearFormats()

$ws.Range("D9").Value = "'0.1231"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'9.52%"
$ws.Range("E9").ClearFormats()
$ws.Range("G9").Value = "'12"
$ws.Range("G9").ClearFormats()

$ws.Range("D10").Value = "'0.1823"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'1.20%"
$ws.Range("E10").ClearFormats()
$ws.Range("G10").Value = "'12"
$ws.Range("G10").ClearFormats()

$ws.Range("D11").Value = "'0.09249"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'1.16%"
$ws.Range("E11").ClearFormats()
$ws.Range("G11").Value = "'12"
$ws.Range("G11").ClearFormats()

$ws.Range("D12").Value = "'0.04192"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-0.66%"
$ws.Range("E12").ClearFormats()
$ws.Range("G12").Value = "'12"
$ws.Range("G12").ClearFormats()

$ws.Range("D13").Value = "'0.1051"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.03%"
$ws.Range("E13").ClearFormats()
$ws.Range("G13").Value = "'12"
$ws.Range("G13").ClearFormats()

$ws.Range("D14").Value = "'0.001253"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'-0.02%"
$ws.Range("E14").ClearFormats()
$ws.Range("G14").Value = "'12"
$ws.Range("G14").ClearFormats()

$ws.Range("D15").Value = "'0.005879"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'1.84%"
$ws.Range("E15").ClearFormats()
$ws.Range("G15").Value = "'12"
$ws.Range("G15").ClearFormats()

$ws.Range("D16").Value = "'3.352"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.40%"
$ws.Range("E16").ClearFormats()
$ws.Range("G16").Value = "'12"
$ws.Range("G16").ClearFormats()

$ws.Range("D17").Value = "'4.319"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'1.83%"
$ws.Range("E17").ClearFormats()
$ws.Range("G17").Value = "'12"
$ws.Range("G17").ClearFormats()

$ws.Range("G18").Value = "'12"
$ws.Range("G18").ClearFormats()

$ws.Range("D19").Value = "'7.583"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'14.40%"
$ws.Range("E19").ClearFormats()
$ws.Range("G19").Value = "'12"
$ws.Range("G19").ClearFormats()

$ws.Range("E20").Value = "'2.86%"
$ws.Range("E20").ClearFormats()
$ws.Range("G20").Value = "'12"
$ws.Range("G20").ClearFormats()

$ws.Range("E21").Value = "'4.38%"
$ws.Range("E21").ClearFormats()
$ws.Range("G21").Value = "'12"
$ws.Range("G21").ClearFormats()

$ws.Range("D22").Value = "'0.04024"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'-0.36%"
$ws.Range("E22").ClearFormats()
$ws.Range("G22").Value = "'12"
$ws.Range("G22").ClearFormats()

$ws.Range("E23").Value = "'1.99%"
$ws.Range("E23").ClearFormats()
$ws.Range("G23").Value = "'12"
$ws.Range("G23").ClearFormats()

$ws.Range("D24").Value = "'0.004089"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'-0.12%"
$ws.Range("E24").ClearFormats()
$ws.Range("G24").Value = "'12"
$ws.Range("G24").ClearFormats()

$ws.Range("D25").Value = "'0.0001302"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'0.15%"
$ws.Range("E25").ClearFormats()
$ws.Range("G25").Value = "'12"
$ws.Range("G25").ClearFormats()

$ws.Range("G26").Value = "'12"
$ws.Range("G26").ClearFormats()

$ws.Range("G27").Value = "'12"
$ws.Range("G27").ClearFormats()

$ws.Range("G28").Value = "'12"
$ws.Range("G28").ClearFormats()

$ws.Range("G29").Value = "'12"
$ws.Range("G29").ClearFormats()

$ws.Range("G30").Value = "'12"
$ws.Range("G30").ClearFormats()

$ws.Range("G31").Value = "'12"
$ws.Range("G31").ClearFormats()

$ws.Range("G32").Value = "'12"
$ws.Range("G32").ClearFormats()

$ws.Range("G33").Value = "'12"
$ws.Range("G33").ClearFormats()

$ws.Range("G34").Value = "'12"
$ws.Range("G34").ClearFormats()

$ws.Range("G35").Value = "'12"
$ws.Range("G35").ClearFormats()

$ws.Range("G36").Value = "'12"
$ws.Range("G36").ClearFormats()

$ws.Range("G37").Value = "'12"
$ws.Range("G37").ClearFormats()

$ws.Range("D38").Value = "'0.02531"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'5.19%"
$ws.Range("E38").ClearFormats()
$ws.Range("G38").Value = "'12"
$ws.Range("G38").ClearFormats()

$ws.Range("D39").Value = "'0.05333"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'1.72%"
$ws.Range("E39").ClearFormats()
$ws.Range("G39").Value = "'12"
$ws.Range("G39").ClearFormats()

$ws.Range("D40").Value = "'0.007851"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'0.24%"
$ws.Range("E40").ClearFormats()
$ws.Range("G40").Value = "'12"
$ws.Range("G40").ClearFormats()

$ws.Range("D41").Value = "'0.1312"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'0.76%"
$ws.Range("E41").ClearFormats()
$ws.Range("G41").Value = "'12"
$ws.Range("G41").ClearFormats()

$ws.Range("D42").Value = "'0.006652"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-5.57%"
$ws.Range("E42").ClearFormats()
$ws.Range("G42").Value = "'12"
$ws.Range("G42").ClearFormats()

$ws.Range("D43").Value = "'0.001853"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-1.98%"
$ws.Range("E43").ClearFormats()
$ws.Range("G43").Value = "'12"
$ws.Range("G43").ClearFormats()

$ws.Range("D44").Value = "'0.007392"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-4.24%"
$ws.Range("E44").ClearFormats()
$ws.Range("G44").Value = "'12"
$ws.Range("G44").ClearFormats()

$ws.Range("D45").Value = "'0.3073"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'-8.38%"
$ws.Range("E45").ClearFormats()
$ws.Range("G45").Value = "'12"
$ws.Range("G45").ClearFormats()

$ws.Range("D46").Value = "'0.00006772"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'-1.86%"
$ws.Range("E46").ClearFormats()
$ws.Range("G46").Value = "'12"
$ws.Range("G46").ClearFormats()

$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'0.29%"
$ws.Range("E47").ClearFormats()
$ws.Range("G47").Value = "'12"
$ws.Range("G47").ClearFormats()

$ws.Range("D48").Value = "'0.2105"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'283.11%"
$ws.Range("E48").ClearFormats()
$ws.Range("G48").Value = "'12"
$ws.Range("G48").ClearFormats()

$ws.Range("G49").Value = "'12"
$ws.Range("G49").ClearFormats()

$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'0.29%"
$ws.Range("E50").ClearFormats()
$ws.Range("G50").Value = "'12"
$ws.Range("G50").ClearFormats()

$ws.Range("D51").Value = "'0.0002003"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'0.29%"
$ws.Range("E51").ClearFormats()
$ws.Range("G51").Value = "'12"
$ws.Range("G51").ClearFormats()
